$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 256.08334
$ws.Range("I17").Value = 110
$ws.Range("J17").Value = 258.55933
$ws.Range("K17").Value = 330
$ws.Range("L17").Value = 775.6779899999999
$ws.Range("M17").Value = -162
$ws.Range("N17").Value = -1111.67799

$ws.Range("H32").Value = 2479.5
$ws.Range("I32").Value = 3898.3333
$ws.Range("J32").Value = 1628.2
$ws.Range("K32").Value = 3898.3333
$ws.Range("L32").Value = 1628.2
$ws.Range("M32").Value = -3572.3333
$ws.Range("N32").Value = -2280.2

$ws.Range("H43").Value = 4112.75
$ws.Range("J43").Value = 4112.75
$ws.Range("L43").Value = 4112.75
$ws.Range("N43").Value = -4250.75

$ws.Range("H51").Value = 2742.762
$ws.Range("I51").Value = 3099.4
$ws.Range("J51").Value = 2631.3125
$ws.Range("K51").Value = 3099.4
$ws.Range("L51").Value = 2631.3125
$ws.Range("M51").Value = -2615.4
$ws.Range("N51").Value = -3599.3125

$ws.Range("H125").Value = 2001.1538
$ws.Range("I125").Value = 2566.2
$ws.Range("J125").Value = 1648
$ws.Range("K125").Value = 23095.8
$ws.Range("L125").Value = 14832
$ws.Range("M125").Value = -20635.8
$ws.Range("N125").Value = -19752

$ws.Range("H132").Value = 273117.12
$ws.Range("I132").Value = 306102.7
$ws.Range("J132").Value = 986.25
$ws.Range("K132").Value = 918308.1000000001
$ws.Range("L132").Value = 2958.75
$ws.Range("M132").Value = -915778.1000000001
$ws.Range("N132").Value = -8018.75

$ws.Range("H137").Value = 43479964
$ws.Range("I137").Value = 1644.579
$ws.Range("J137").Value = 250001970
$ws.Range("K137").Value = 4933.737
$ws.Range("L137").Value = 750005910
$ws.Range("M137").Value = -2383.737
$ws.Range("N137").Value = -750011010

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2101907
$ws.Range("I2").Value = 863.5
$ws.Range("J2").Value = 7354516
$ws.Range("K2").Value = 863.5
$ws.Range("L2").Value = 7354516
$ws.Range("M2").Value = -750.5
$ws.Range("N2").Value = -7354742

$ws.Range("H32").Value = 5762.259
$ws.Range("I32").Value = 6355.8667
$ws.Range("J32").Value = 3707.4614
$ws.Range("K32").Value = 6355.8667
$ws.Range("L32").Value = 3707.4614
$ws.Range("M32").Value = -6068.8667
$ws.Range("N32").Value = -4281.4614

$ws.Range("H45").Value = 820.3
$ws.Range("I45").Value = 820.3
$ws.Range("K45").Value = 820.3
$ws.Range("M45").Value = -443.3

$ws.Range("H116").Value = 2101907
$ws.Range("I116").Value = 863.5
$ws.Range("J116").Value = 7354516
$ws.Range("K116").Value = 863.5
$ws.Range("L116").Value = 7354516
$ws.Range("M116").Value = 1430.5
$ws.Range("N116").Value = -7359104

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2101907
$ws.Range("I3").Value = 863.5
$ws.Range("J3").Value = 7354516
$ws.Range("K3").Value = 863.5
$ws.Range("L3").Value = 7354516
$ws.Range("M3").Value = -749.5
$ws.Range("N3").Value = -7354744

$ws.Range("H134").Value = 84595.25
$ws.Range("I134").Value = 92176.63
$ws.Range("J134").Value = 1200
$ws.Range("K134").Value = 276529.89
$ws.Range("L134").Value = 3600
$ws.Range("M134").Value = -273994.89
$ws.Range("N134").Value = -8670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1944
$ws.Range("I31").Value = 2021.7142
$ws.Range("K31").Value = 2021.7142
$ws.Range("M31").Value = -1726.7142

$ws.Range("H34").Value = 1944
$ws.Range("I34").Value = 2021.7142
$ws.Range("K34").Value = 2021.7142
$ws.Range("M34").Value = -1819.7142

$ws.Range("H122").Value = 727.2353000000001
$ws.Range("I122").Value = 713.75
$ws.Range("J122").Value = 739.2222
$ws.Range("K122").Value = 2141.25
$ws.Range("L122").Value = 2217.6666
$ws.Range("M122").Value = 308.75
$ws.Range("N122").Value = -7117.6666

$ws.Range("H132").Value = 4150.45
$ws.Range("I132").Value = 3467.8
$ws.Range("J132").Value = 6198.4
$ws.Range("K132").Value = 10403.4
$ws.Range("L132").Value = 18595.2
$ws.Range("M132").Value = -7873.400000000001
$ws.Range("N132").Value = -23655.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 83334190
$ws.Range("J34").Value = 111112104
$ws.Range("L34").Value = 333336312
$ws.Range("N34").Value = -333336480

$ws.Range("H108").Value = 1000
$ws.Range("I108").Value = 1000
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 3000
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -120
$ws.Range("N108").ClearContents()

$ws.Range("H131").Value = 3135.0227
$ws.Range("J131").Value = 2329.4866
$ws.Range("L131").Value = 6988.459800000001
$ws.Range("N131").Value = -17068.4598

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 4407.5
$ws.Range("I55").Value = 4065
$ws.Range("J55").Value = 4750
$ws.Range("K55").Value = 4065
$ws.Range("L55").Value = 4750
$ws.Range("M55").Value = -3738
$ws.Range("N55").Value = -5404

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 393.8
$ws.Range("I22").Value = 356.33334
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 356.33334
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = -61.33334000000002
$ws.Range("N22").Value = -1040

$ws.Range("H27").Value = 393.8
$ws.Range("I27").Value = 356.33334
$ws.Range("J27").Value = 450
$ws.Range("K27").Value = 356.33334
$ws.Range("L27").Value = 450
$ws.Range("M27").Value = -249.33334
$ws.Range("N27").Value = -664

$ws.Range("H40").Value = 2508.3333
$ws.Range("I40").Value = 2230.4443
$ws.Range("J40").Value = 2925.1667
$ws.Range("K40").Value = 2230.4443
$ws.Range("L40").Value = 2925.1667
$ws.Range("M40").Value = -2094.4443
$ws.Range("N40").Value = -3197.1667

$ws.Range("H127").Value = 48097.5
$ws.Range("J127").Value = 48097.5
$ws.Range("L127").Value = 48097.5
$ws.Range("N127").Value = -58017.5

$ws.Range("H136").Value = 900.6415
$ws.Range("I136").Value = 762.1951
$ws.Range("J136").Value = 1373.6666
$ws.Range("K136").Value = 2286.5853
$ws.Range("L136").Value = 4120.9998
$ws.Range("M136").Value = 263.4146999999998
$ws.Range("N136").Value = -9220.9998

$ws.Range("H141").Value = 50000
$ws.Range("I141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("K141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("M141").Value = -44820
$ws.Range("N141").Value = -60360
